# The catalog/tag parser now computes sale prices itself, so the manually
# entered "price" column and the sample/demo data rows it was maintained in
# are no longer needed. Remove them from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two demo/sample data rows (row 2 and row 3), leaving only the
# header row.
$ws.Rows("2:3").Delete()

# Drop column B ("price") entirely; the remaining quantity columns shift
# left to take its place.
$ws.Columns("B").Delete()
